$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2340.4746
$ws.Range("J17").Value = 2340.4746
$ws.Range("L17").Value = 7021.4238
$ws.Range("N17").Value = -7357.4238

$ws.Range("H21").Value = 70014.25
$ws.Range("J21").Value = 70000
$ws.Range("L21").Value = 70000
$ws.Range("N21").Value = -70936

$ws.Range("H23").Value = 70014.25
$ws.Range("J23").Value = 70000
$ws.Range("L23").Value = 70000
$ws.Range("N23").Value = -70468

$ws.Range("H98").Value = 99532.125
$ws.Range("I98").Value = 1993.3334
$ws.Range("K98").Value = 1993.3334
$ws.Range("M98").Value = -495.3334

$ws.Range("H101").Value = 734
$ws.Range("I101").Value = 734
$ws.Range("K101").Value = 2202
$ws.Range("M101").Value = -580

$ws.Range("H122").Value = 99532.125
$ws.Range("I122").Value = 1993.3334
$ws.Range("K122").Value = 5980.0002
$ws.Range("M122").Value = -3530.0002

$ws.Range("H132").Value = 25478.85
$ws.Range("I132").Value = 3488.9062
$ws.Range("K132").Value = 10466.7186
$ws.Range("M132").Value = -7936.7186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 429991.53
$ws.Range("I45").Value = 601263.5600000001
$ws.Range("J45").Value = 1811.5
$ws.Range("K45").Value = 601263.5600000001
$ws.Range("L45").Value = 1811.5
$ws.Range("M45").Value = -600886.5600000001
$ws.Range("N45").Value = -2565.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3388.239
$ws.Range("I134").Value = 3063
$ws.Range("J134").Value = 3561.7
$ws.Range("K134").Value = 9189
$ws.Range("L134").Value = 10685.1
$ws.Range("M134").Value = -6654
$ws.Range("N134").Value = -15755.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 1466
$ws.Range("I36").Value = 1466
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1466
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1078
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 1466
$ws.Range("I40").Value = 1466
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1466
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1306
$ws.Range("N40").ClearContents()

$ws.Range("H58").Value = 1671.9131
$ws.Range("I58").Value = 1093.2354
$ws.Range("J58").Value = 3311.5
$ws.Range("K58").Value = 1093.2354
$ws.Range("L58").Value = 3311.5
$ws.Range("M58").Value = -890.2354
$ws.Range("N58").Value = -3717.5

$ws.Range("H69").Value = 18128.2
$ws.Range("I69").Value = 10641
$ws.Range("K69").Value = 10641
$ws.Range("M69").Value = -9892

$ws.Range("H72").Value = 18128.2
$ws.Range("I72").Value = 10641
$ws.Range("K72").Value = 31923
$ws.Range("M72").Value = -28179

$ws.Range("H99").Value = 2617.9
$ws.Range("I99").Value = 2294.7
$ws.Range("J99").Value = 2779.5
$ws.Range("K99").Value = 2294.7
$ws.Range("L99").Value = 2779.5
$ws.Range("M99").Value = -796.6999999999998
$ws.Range("N99").Value = -5775.5

$ws.Range("H126").Value = 2617.9
$ws.Range("I126").Value = 2294.7
$ws.Range("J126").Value = 2779.5
$ws.Range("K126").Value = 6884.099999999999
$ws.Range("L126").Value = 8338.5
$ws.Range("M126").Value = -4414.099999999999
$ws.Range("N126").Value = -13278.5

$ws.Range("H136").Value = 1671.9131
$ws.Range("I136").Value = 1093.2354
$ws.Range("J136").Value = 3311.5
$ws.Range("K136").Value = 3279.7062
$ws.Range("L136").Value = 9934.5
$ws.Range("M136").Value = -729.7062000000001
$ws.Range("N136").Value = -15034.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 200
$ws.Range("I24").Value = 200
$ws.Range("K24").Value = 600
$ws.Range("M24").Value = -370

$ws.Range("H34").Value = 1612.375
$ws.Range("J34").Value = 1799.8572
$ws.Range("L34").Value = 5399.571599999999
$ws.Range("N34").Value = -5567.571599999999

$ws.Range("H39").Value = 2820
$ws.Range("J39").Value = 2820
$ws.Range("L39").Value = 8460
$ws.Range("N39").Value = -9048

$ws.Range("H40").Value = 4180.5386
$ws.Range("I40").Value = 5660.4443
$ws.Range("J40").Value = 850.75
$ws.Range("K40").Value = 22641.7772
$ws.Range("L40").Value = 3403
$ws.Range("M40").Value = -22572.7772
$ws.Range("N40").Value = -3541

$ws.Range("H98").Value = 1700
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1700
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 5100
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -8096

$ws.Range("H113").Value = 6649.353
$ws.Range("J113").Value = 998.8
$ws.Range("L113").Value = 2996.4
$ws.Range("N113").Value = -7336.4

$ws.Range("H116").Value = 2000
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H129").Value = 108412.54
$ws.Range("I129").Value = 333997.78
$ws.Range("J129").Value = 1556.3684
$ws.Range("K129").Value = 1001993.34
$ws.Range("L129").Value = 4669.1052
$ws.Range("M129").Value = -996993.3400000001
$ws.Range("N129").Value = -14669.1052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 8106.8335
$ws.Range("I31").Value = 1728.2
$ws.Range("K31").Value = 1728.2
$ws.Range("M31").Value = -1436.2

$ws.Range("H37").Value = 8106.8335
$ws.Range("I37").Value = 1728.2
$ws.Range("K37").Value = 1728.2
$ws.Range("M37").Value = -1451.2

$ws.Range("H126").Value = 76926040
$ws.Range("I126").Value = 250003400
$ws.Range("J126").Value = 2767.111
$ws.Range("K126").Value = 750010200
$ws.Range("L126").Value = 8301.332999999999
$ws.Range("M126").Value = -750007730
$ws.Range("N126").Value = -13241.333

$ws.Range("H132").Value = 35716910
$ws.Range("I132").Value = 58825268
$ws.Range("J132").Value = 4000.2727
$ws.Range("K132").Value = 176475804
$ws.Range("L132").Value = 12000.8181
$ws.Range("M132").Value = -176473274
$ws.Range("N132").Value = -17060.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 79423.38
$ws.Range("I122").Value = 102400.4
$ws.Range("J122").Value = 2833.3333
$ws.Range("K122").Value = 307201.2
$ws.Range("L122").Value = 8499.999899999999
$ws.Range("M122").Value = -304751.2
$ws.Range("N122").Value = -13399.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 19112.5
$ws.Range("J69").Value = 19112.5
$ws.Range("L69").Value = 19112.5
$ws.Range("N69").Value = -20610.5

$ws.Range("H72").Value = 19112.5
$ws.Range("J72").Value = 19112.5
$ws.Range("L72").Value = 57337.5
$ws.Range("N72").Value = -64825.5

$ws.Range("H122").Value = 9524876
$ws.Range("I122").Value = 28571428
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 85714284
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -85711834
$ws.Range("N122").Value = -9700
